$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.812.37'
$ws.Range("E2").Value = '  +0.13%  '

$ws.Range("D3").Value = '2.538.23'
$ws.Range("E3").Value = '  -0.41%  '

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = '''304.79'
$ws.Range("E5").Value = '  +1.64%  '

$ws.Range("D6").Value = '''97.41'
$ws.Range("E6").Value = '  +5.24%  '

$ws.Range("D7").Value = '''0.576'
$ws.Range("E7").Value = '  +0.41%  '

$ws.Range("E8").Value = '  +0.09%  '

$ws.Range("E9").Value = '  -1.00%  '

$ws.Range("D10").Value = '''36.44'
$ws.Range("E10").Value = '  +1.14%  '

$ws.Range("E11").Value = '  +2.11%  '

$ws.Range("E12").Value = '  +0.61%  '

$ws.Range("D13").Value = '''7.55'
$ws.Range("E13").Value = '  -2.18%  '

$ws.Range("D14").Value = '2.930.35'
$ws.Range("E14").Value = '  -0.17%  '

$ws.Range("D15").Value = '2.580.74'
$ws.Range("E15").Value = '  +2.45%  '

$ws.Range("D16").Value = '''14.99'
$ws.Range("E16").Value = '  +5.42%  '

$ws.Range("D17").Value = '''0.863'
$ws.Range("E17").Value = '  -1.54%  '

$ws.Range("D18").Value = '42.816.98'
$ws.Range("E18").Value = '  +0.14%  '

$ws.Range("D19").Value = '''13.23'
$ws.Range("E19").Value = '  +2.55%  '

$ws.Range("E20").Value = '  +0.14%  '

$ws.Range("D21").Value = '''6.55'
$ws.Range("E21").Value = '  -0.48%  '

$ws.Range("D22").Value = '''71.60'
$ws.Range("E22").Value = '  -0.28%  '

$ws.Range("D23").Value = '''253.56'
$ws.Range("E23").Value = '  -0.70%  '

$ws.Range("E24").Value = '  +0.10%  '

$ws.Range("E25").Value = '  -2.64%  '

$ws.Range("D26").Value = '''27.84'
$ws.Range("E26").Value = '  -4.82%  '

$ws.Range("D27").Value = '''0.998'
$ws.Range("E27").Value = '  -0.29%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '''2.27'
$ws.Range("E28").Value = '  +7.49%  '

$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").Value = '''10.16'
$ws.Range("E29").Value = '  +0.34%  '

$ws.Range("D30").Value = '''37.86'
$ws.Range("E30").Value = '  +2.26%  '

$ws.Range("D31").Value = '''6.15'
$ws.Range("E31").Value = '  +2.16%  '

$ws.Range("D32").Value = '''157.17'
$ws.Range("E32").Value = '  +3.05%  '

$ws.Range("D33").Value = '''19.47'
$ws.Range("E33").Value = '  +13.45%  '

$ws.Range("D34").Value = '''2.13'
$ws.Range("E34").Value = '  -1.87%  '

$ws.Range("D35").Value = '''3.30'
$ws.Range("E35").Value = '  -2.28%  '

$ws.Range("D36").Value = '''0.0793'
$ws.Range("E36").Value = '  -0.14%  '

$ws.Range("E37").Value = '  -4.66%  '

$ws.Range("D38").Value = '''0.114'
$ws.Range("E38").Value = '  -0.14%  '

$ws.Range("D39").Value = '''25.05'
$ws.Range("E39").Value = '  +3.58%  '

$ws.Range("E40").Value = '  +0.29%  '

$ws.Range("D41").Value = '''2.16'
$ws.Range("E41").Value = '  +30.65%  '

$ws.Range("D42").Value = '''3.40'
$ws.Range("E42").Value = '  -0.78%  '

$ws.Range("D43").Value = '''3.85'
$ws.Range("E43").Value = '  -0.59%  '

$ws.Range("D44").Value = '2.093.00'
$ws.Range("E44").Value = '  +0.36%  '

$ws.Range("E45").Value = '  -2.33%  '

$ws.Range("E46").Value = '  +0.08%  '

$ws.Range("D47").Value = '''86.20'
$ws.Range("E47").Value = '  +1.96%  '

$ws.Range("E48").Value = '  -1.65%  '

$ws.Range("D49").Value = '2.786.34'
$ws.Range("E49").Value = '  -0.09%  '

$ws.Range("D50").Value = '''73.48'
$ws.Range("E50").Value = '  +6.12%  '

$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").Value = '''0.191'
$ws.Range("E51").Value = '  +1.30%  '
